$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns involved in the weekly rotation of rows 2, 3 and 5
$cols = @("D", "J", "K", "L", "M", "N", "O", "P", "Q")

# Capture the current ("before") values for each row/column that changes
$row2 = @{}
$row3 = @{}
$row5 = @{}
foreach ($c in $cols) {
    $row2[$c] = $ws.Range("${c}2").Value()
    $row3[$c] = $ws.Range("${c}3").Value()
    $row5[$c] = $ws.Range("${c}5").Value()
}

# Apply the cyclic rotation: old row3 -> row2, old row5 -> row3, old row2 -> row5
foreach ($c in $cols) {
    $ws.Range("${c}2").Value = $row3[$c]
    $ws.Range("${c}3").Value = $row5[$c]
    $ws.Range("${c}5").Value = $row2[$c]
}
